$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate B:AC values among swapped/shuffled match rows ---
# Group: [98, 99]
$r98 = $ws.Range("B98:AC98")
$r99 = $ws.Range("B99:AC99")
$v98 = $r98.Value2
$v99 = $r99.Value2
$r98.Value2 = $v99
$r99.Value2 = $v98

# Group: [105, 106]
$r105 = $ws.Range("B105:AC105")
$r106 = $ws.Range("B106:AC106")
$v105 = $r105.Value2
$v106 = $r106.Value2
$r105.Value2 = $v106
$r106.Value2 = $v105

# Group: [139, 140]
$r139 = $ws.Range("B139:AC139")
$r140 = $ws.Range("B140:AC140")
$v139 = $r139.Value2
$v140 = $r140.Value2
$r139.Value2 = $v140
$r140.Value2 = $v139

# Group: [142, 143, 144]
$r142 = $ws.Range("B142:AC142")
$r143 = $ws.Range("B143:AC143")
$r144 = $ws.Range("B144:AC144")
$v142 = $r142.Value2
$v143 = $r143.Value2
$v144 = $r144.Value2
$r142.Value2 = $v143
$r143.Value2 = $v144
$r144.Value2 = $v142

# Group: [149, 150]
$r149 = $ws.Range("B149:AC149")
$r150 = $ws.Range("B150:AC150")
$v149 = $r149.Value2
$v150 = $r150.Value2
$r149.Value2 = $v150
$r150.Value2 = $v149

# Group: [155, 156, 157]
$r155 = $ws.Range("B155:AC155")
$r156 = $ws.Range("B156:AC156")
$r157 = $ws.Range("B157:AC157")
$v155 = $r155.Value2
$v156 = $r156.Value2
$v157 = $r157.Value2
$r155.Value2 = $v156
$r156.Value2 = $v157
$r157.Value2 = $v155

# Group: [162, 163, 164]
$r162 = $ws.Range("B162:AC162")
$r163 = $ws.Range("B163:AC163")
$r164 = $ws.Range("B164:AC164")
$v162 = $r162.Value2
$v163 = $r163.Value2
$v164 = $r164.Value2
$r162.Value2 = $v163
$r163.Value2 = $v164
$r164.Value2 = $v162

# --- Append new rows 165-170 ---
# Row 165
$src = $ws.Range("A164:AC164")
$dst = $ws.Range("A165:AC165")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Cells.Item(165,1).Value2 = 163
$ws.Cells.Item(165,2).Value2 = 7799445
$ws.Cells.Item(165,3).Value2 = "South Africa Premier"
$ws.Cells.Item(165,4).Value2 = "South Africa Premier"
$ws.Cells.Item(165,5).Value2 = 45388.41666666666
$ws.Cells.Item(165,6).Value2 = "Chippa United"
$ws.Cells.Item(165,7).Value2 = "Kaizer Chiefs"
$ws.Cells.Item(165,11).Value2 = 3.3
$ws.Cells.Item(165,12).Value2 = 2.875
$ws.Cells.Item(165,13).Value2 = 2.375
$ws.Cells.Item(165,14).Value2 = 3.6
$ws.Cells.Item(165,15).Value2 = 2.875
$ws.Cells.Item(165,16).Value2 = 2.2
$ws.Cells.Item(165,17).Value2 = 0.25
$ws.Cells.Item(165,18).Value2 = 1.95
$ws.Cells.Item(165,19).Value2 = 1.85
$ws.Cells.Item(165,20).Value2 = 2
$ws.Cells.Item(165,21).Value2 = 1.975
$ws.Cells.Item(165,22).Value2 = 1.825
$ws.Cells.Item(165,23).Value2 = 0
$ws.Cells.Item(165,24).Value2 = 0
$ws.Cells.Item(165,25).Value2 = 0
$ws.Cells.Item(165,26).Value2 = 0
$ws.Cells.Item(165,27).Value2 = 0

# Row 166
$src = $ws.Range("A165:AC165")
$dst = $ws.Range("A166:AC166")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Cells.Item(166,1).Value2 = 164
$ws.Cells.Item(166,2).Value2 = 8052475
$ws.Cells.Item(166,3).Value2 = "South Africa Premier"
$ws.Cells.Item(166,4).Value2 = "South Africa Premier"
$ws.Cells.Item(166,5).Value2 = 45388.5
$ws.Cells.Item(166,6).Value2 = "Supersport United"
$ws.Cells.Item(166,7).Value2 = "TS Galaxy"
$ws.Cells.Item(166,11).Value2 = 1.3
$ws.Cells.Item(166,12).Value2 = 4.75
$ws.Cells.Item(166,13).Value2 = 11
$ws.Cells.Item(166,14).Value2 = 1.571
$ws.Cells.Item(166,15).Value2 = 3.4
$ws.Cells.Item(166,16).Value2 = 6
$ws.Cells.Item(166,17).Value2 = -0.75
$ws.Cells.Item(166,18).Value2 = 1.775
$ws.Cells.Item(166,19).Value2 = 2.025
$ws.Cells.Item(166,20).Value2 = 2
$ws.Cells.Item(166,21).Value2 = 1.8
$ws.Cells.Item(166,22).Value2 = 2
$ws.Cells.Item(166,23).Value2 = 0
$ws.Cells.Item(166,24).Value2 = 0
$ws.Cells.Item(166,25).Value2 = 0
$ws.Cells.Item(166,26).Value2 = 0
$ws.Cells.Item(166,27).Value2 = 0

# Row 167
$src = $ws.Range("A166:AC166")
$dst = $ws.Range("A167:AC167")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Cells.Item(167,1).Value2 = 165
$ws.Cells.Item(167,2).Value2 = 7799663
$ws.Cells.Item(167,3).Value2 = "South Africa Premier"
$ws.Cells.Item(167,4).Value2 = "South Africa Premier"
$ws.Cells.Item(167,5).Value2 = 45388.5
$ws.Cells.Item(167,6).Value2 = "Orlando Pirates"
$ws.Cells.Item(167,7).Value2 = "Golden Arrows"
$ws.Cells.Item(167,11).Value2 = 1.5
$ws.Cells.Item(167,12).Value2 = 3.75
$ws.Cells.Item(167,13).Value2 = 6
$ws.Cells.Item(167,14).Value2 = 1.533
$ws.Cells.Item(167,15).Value2 = 3.75
$ws.Cells.Item(167,16).Value2 = 5.75
$ws.Cells.Item(167,17).Value2 = -1
$ws.Cells.Item(167,18).Value2 = 1.95
$ws.Cells.Item(167,19).Value2 = 1.85
$ws.Cells.Item(167,20).Value2 = 2.25
$ws.Cells.Item(167,21).Value2 = 1.775
$ws.Cells.Item(167,22).Value2 = 2.025
$ws.Cells.Item(167,23).Value2 = 0
$ws.Cells.Item(167,24).Value2 = 0
$ws.Cells.Item(167,25).Value2 = 0
$ws.Cells.Item(167,26).Value2 = 0
$ws.Cells.Item(167,27).Value2 = 0

# Row 168
$src = $ws.Range("A167:AC167")
$dst = $ws.Range("A168:AC168")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Cells.Item(168,1).Value2 = 166
$ws.Cells.Item(168,2).Value2 = 7799661
$ws.Cells.Item(168,3).Value2 = "South Africa Premier"
$ws.Cells.Item(168,4).Value2 = "South Africa Premier"
$ws.Cells.Item(168,5).Value2 = 45388.60416666666
$ws.Cells.Item(168,6).Value2 = "Stellenbosch FC"
$ws.Cells.Item(168,7).Value2 = "Sekhukhune United FC"
$ws.Cells.Item(168,11).Value2 = 2.05
$ws.Cells.Item(168,12).Value2 = 3.1
$ws.Cells.Item(168,13).Value2 = 4
$ws.Cells.Item(168,14).Value2 = 2.05
$ws.Cells.Item(168,15).Value2 = 3.1
$ws.Cells.Item(168,16).Value2 = 4
$ws.Cells.Item(168,17).Value2 = -0.5
$ws.Cells.Item(168,18).Value2 = 2.025
$ws.Cells.Item(168,19).Value2 = 1.775
$ws.Cells.Item(168,20).Value2 = 2
$ws.Cells.Item(168,21).Value2 = 1.95
$ws.Cells.Item(168,22).Value2 = 1.85
$ws.Cells.Item(168,23).Value2 = 0
$ws.Cells.Item(168,24).Value2 = 0
$ws.Cells.Item(168,25).Value2 = 0
$ws.Cells.Item(168,26).Value2 = 0
$ws.Cells.Item(168,27).Value2 = 0

# Row 169
$src = $ws.Range("A168:AC168")
$dst = $ws.Range("A169:AC169")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Cells.Item(169,1).Value2 = 167
$ws.Cells.Item(169,2).Value2 = 7802375
$ws.Cells.Item(169,3).Value2 = "South Africa Premier"
$ws.Cells.Item(169,4).Value2 = "South Africa Premier"
$ws.Cells.Item(169,5).Value2 = 45389.41666666666
$ws.Cells.Item(169,6).Value2 = "Royal AM FC"
$ws.Cells.Item(169,7).Value2 = "Polokwane City"
$ws.Cells.Item(169,11).Value2 = 2.4
$ws.Cells.Item(169,12).Value2 = 2.875
$ws.Cells.Item(169,13).Value2 = 3
$ws.Cells.Item(169,14).Value2 = 2.8
$ws.Cells.Item(169,15).Value2 = 2.875
$ws.Cells.Item(169,16).Value2 = 2.6
$ws.Cells.Item(169,17).Value2 = 0
$ws.Cells.Item(169,18).Value2 = 1.975
$ws.Cells.Item(169,19).Value2 = 1.825
$ws.Cells.Item(169,20).Value2 = 2
$ws.Cells.Item(169,21).Value2 = 1.925
$ws.Cells.Item(169,22).Value2 = 1.875
$ws.Cells.Item(169,23).Value2 = 0
$ws.Cells.Item(169,24).Value2 = 0
$ws.Cells.Item(169,25).Value2 = 0
$ws.Cells.Item(169,26).Value2 = 0
$ws.Cells.Item(169,27).Value2 = 0

# Row 170
$src = $ws.Range("A169:AC169")
$dst = $ws.Range("A170:AC170")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Cells.Item(170,1).Value2 = 168
$ws.Cells.Item(170,2).Value2 = 8052476
$ws.Cells.Item(170,3).Value2 = "South Africa Premier"
$ws.Cells.Item(170,4).Value2 = "South Africa Premier"
$ws.Cells.Item(170,5).Value2 = 45389.5
$ws.Cells.Item(170,6).Value2 = "Cape Town City"
$ws.Cells.Item(170,7).Value2 = "Amazulu"
$ws.Cells.Item(170,11).Value2 = 2.05
$ws.Cells.Item(170,12).Value2 = 3.4
$ws.Cells.Item(170,13).Value2 = 3.3
$ws.Cells.Item(170,14).Value2 = 2.05
$ws.Cells.Item(170,15).Value2 = 3.4
$ws.Cells.Item(170,16).Value2 = 3.3
$ws.Cells.Item(170,17).Value2 = -0.25
$ws.Cells.Item(170,18).Value2 = 1.8
$ws.Cells.Item(170,19).Value2 = 2
$ws.Cells.Item(170,20).Value2 = 2
$ws.Cells.Item(170,21).Value2 = 1.95
$ws.Cells.Item(170,22).Value2 = 1.85
$ws.Cells.Item(170,23).Value2 = 0
$ws.Cells.Item(170,24).Value2 = 0
$ws.Cells.Item(170,25).Value2 = 0
$ws.Cells.Item(170,26).Value2 = 0
$ws.Cells.Item(170,27).Value2 = 0

